# "Add files via upload" — the two SharePoint-metadata customXml parts got
# swapped when the deck was re-uploaded:
#   - customXml/item1.xml used to hold the FormTemplates / "mso-contentType"
#     document-library form stub (namespace
#     http://schemas.microsoft.com/sharepoint/v3/contenttype/forms)
#   - customXml/item2.xml used to hold the p:properties/documentManagement
#     SharePoint column data (namespace
#     http://schemas.microsoft.com/office/2006/metadata/properties)
# After the upload, item1.xml carries the p:properties payload and item2.xml
# carries the FormTemplates payload — i.e. the XML contents of the two parts
# were exchanged. (Each part's paired itemProps*.xml datastore item, wired
# through customXml/_rels/item*.xml.rels, is maintained automatically by the
# CustomXMLParts collection and follows its owning part.)
#
# Reproduce that swap through the CustomXMLParts COM surface: look each part
# up by the namespace that identifies its content (not by its current
# positional index, since that's exactly what got shuffled), capture the two
# XML payloads, delete both parts, then re-add them with the payloads
# exchanged.

$p = $ppt.ActivePresentation
$parts = $p.CustomXMLParts

$formsNamespace = "http://schemas.microsoft.com/sharepoint/v3/contenttype/forms"
$propsNamespace = "http://schemas.microsoft.com/office/2006/metadata/properties"

$formsPart = $parts.SelectByNamespace($formsNamespace).Item(1)
$propsPart = $parts.SelectByNamespace($propsNamespace).Item(1)

$formsXml = $formsPart.XML
$propsXml = $propsPart.XML

$formsPart.Delete()
$propsPart.Delete()

# Re-add with the XML exchanged between the two parts (item1 <-> item2).
$parts.Add($propsXml) | Out-Null
$parts.Add($formsXml) | Out-Null
